$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.466.96"
$ws.Range("E2").Value = "  -1.35%  "

# Row 3
$ws.Range("D3").Value = "2.454.61"
$ws.Range("E3").Value = "  -0.05%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.73"
$ws.Range("E5").Value = "  +1.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.46"
$ws.Range("E6").Value = "  -1.58%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  -1.08%  "

# Row 9
$ws.Range("D9").Value = "2.451.15"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("E10").Value = "  -3.58%  "

# Row 11
$ws.Range("E11").Value = "  +2.02%  "

# Row 12
$ws.Range("E12").Value = "  -1.19%  "

# Row 13
$ws.Range("E13").Value = "  -3.28%  "

# Row 14
$ws.Range("E14").Value = "  -1.70%  "

# Row 15
$ws.Range("E15").Value = "  -4.04%  "

# Row 16
$ws.Range("D16").Value = "2.913.30"
$ws.Range("E16").Value = "  +0.46%  "

# Row 17
$ws.Range("D17").Value = "62.311.95"
$ws.Range("E17").Value = "  -1.46%  "

# Row 18
$ws.Range("D18").Value = "2.451.76"
$ws.Range("E18").Value = "  -0.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.95"
$ws.Range("E19").Value = "  -3.22%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("E20").Value = "  -2.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.30"
$ws.Range("E21").Value = "  +0.58%  "

# Row 22
$ws.Range("E22").Value = "  -2.21%  "

# Row 23
$ws.Range("E23").Value = "  -3.93%  "

# Row 24
$ws.Range("E24").Value = "  +0.23%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.11"
$ws.Range("E25").Value = "  +0.60%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.48"
$ws.Range("E26").Value = "  +4.90%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "629.89"

# Row 28
$ws.Range("D28").Value = "0.0₃0965"
$ws.Range("E28").Value = "  -6.20%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.13%  "

# Row 31
$ws.Range("E31").Value = "  -4.56%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.05"
$ws.Range("E32").Value = "  -2.73%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.142"
$ws.Range("E33").Value = "  +0.89%  "

# Row 34
$ws.Range("E34").Value = "  -0.74%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.95"
$ws.Range("E35").Value = "  -5.11%  "

# Row 36
$ws.Range("E36").Value = "  -0.04%  "

# Row 37
$ws.Range("E37").Value = "  -6.25%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.379"
$ws.Range("E38").Value = "  -0.39%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.36"
$ws.Range("E39").Value = "  -1.82%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.86"
$ws.Range("E40").Value = "  +1.99%  "

# Row 41
$ws.Range("E41").Value = "  -2.72%  "

# Row 42
$ws.Range("E42").Value = "  -2.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.49"
$ws.Range("E43").Value = "  +1.64%  "

# Row 45
$ws.Range("E45").Value = "  -4.90%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.22"
$ws.Range("E46").Value = "  -3.17%  "

# Row 47
$ws.Range("E47").Value = "  -3.07%  "

# Row 48
$ws.Range("E48").Value = "  -1.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.603"
$ws.Range("E49").Value = "  +0.09%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.75"
$ws.Range("E50").Value = "  -7.21%  "

# Row 51
$ws.Range("D51").Value = "0.0₆0236"
$ws.Range("E51").Value = "  +6.90%  "
